$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.83830466666667
$ws.Range("N2").Value = 41.514914
$ws.Range("O2").Value = 0.1568893148900199
$ws.Range("P2").Value = 0.1568893148900199
$ws.Range("Q2").Value = 2787.584530996266
$ws.Range("R2").Value = 25088.2607789664
$ws.Range("S2").Value = 0.07583245400896756
$ws.Range("T2").Value = 0.07583245400896758

# Row 3
$ws.Range("G3").Value = 201.4397426666667
$ws.Range("H3").Value = 604.3192280000001
$ws.Range("I3").Value = 0.4833500233086392
$ws.Range("J3").Value = 0.4833500233086393
$ws.Range("O3").Value = 0.5509859018285573
$ws.Range("P3").Value = 0.5509859018285573
$ws.Range("Q3").Value = 9789.830351486973
$ws.Range("R3").Value = 88108.47316338276
$ws.Range("S3").Value = 0.2663190484915648
$ws.Range("T3").Value = 0.2663190484915648

# Row 4
$ws.Range("G4").Value = 201.4397426666667
$ws.Range("H4").Value = 604.3192280000001
$ws.Range("I4").Value = 0.4833500233086392
$ws.Range("J4").Value = 0.4833500233086393
$ws.Range("M4").Value = 25.766648
$ws.Range("N4").Value = 77.299944
$ws.Range("O4").Value = 0.2921247832814228
$ws.Range("P4").Value = 0.2921247832814228
$ws.Range("Q4").Value = 5190.426942502582
$ws.Range("R4").Value = 46713.84248252324
$ws.Range("S4").Value = 0.1411985208081069
$ws.Range("T4").Value = 0.1411985208081069

# Row 5
$ws.Range("I5").Value = 0.1569674599353791
$ws.Range("J5").Value = 0.1569674599353792
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.83830466666667
$ws.Range("N5").Value = 41.514914
$ws.Range("O5").Value = 0.1568893148900199
$ws.Range("P5").Value = 0.1568893148900199
$ws.Range("Q5").Value = 905.2654227477685
$ws.Range("R5").Value = 8147.388804729917
$ws.Range("S5").Value = 0.02462651724928828
$ws.Range("T5").Value = 0.02462651724928829

# Row 6
$ws.Range("I6").Value = 0.1569674599353791
$ws.Range("J6").Value = 0.1569674599353792
$ws.Range("O6").Value = 0.5509859018285573
$ws.Range("P6").Value = 0.5509859018285573
$ws.Range("S6").Value = 0.08648685747023281
$ws.Range("T6").Value = 0.08648685747023283

# Row 7
$ws.Range("I7").Value = 0.1569674599353791
$ws.Range("J7").Value = 0.1569674599353792
$ws.Range("M7").Value = 25.766648
$ws.Range("N7").Value = 77.299944
$ws.Range("O7").Value = 0.2921247832814228
$ws.Range("P7").Value = 0.2921247832814228
$ws.Range("Q7").Value = 1685.586208453637
$ws.Range("R7").Value = 15170.27587608273
$ws.Range("S7").Value = 0.04585408521585804
$ws.Range("T7").Value = 0.04585408521585805

# Row 8
$ws.Range("G8").Value = 60.43484133333334
$ws.Range("H8").Value = 181.304524
$ws.Range("I8").Value = 0.1450120099461104
$ws.Range("J8").Value = 0.1450120099461104
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 13.83830466666667
$ws.Range("N8").Value = 41.514914
$ws.Range("O8").Value = 0.1568893148900199
$ws.Range("P8").Value = 0.1568893148900199
$ws.Range("Q8").Value = 836.3157468523264
$ws.Range("R8").Value = 7526.841721670938
$ws.Range("S8").Value = 0.02275083489127
$ws.Range("T8").Value = 0.02275083489127001

# Row 9
$ws.Range("G9").Value = 60.43484133333334
$ws.Range("H9").Value = 181.304524
$ws.Range("I9").Value = 0.1450120099461104
$ws.Range("J9").Value = 0.1450120099461104
$ws.Range("O9").Value = 0.5509859018285573
$ws.Range("P9").Value = 0.5509859018285573
$ws.Range("Q9").Value = 2937.090944121173
$ws.Range("R9").Value = 26433.81849709055
$ws.Range("S9").Value = 0.07989957307612934
$ws.Range("T9").Value = 0.07989957307612935

# Row 10
$ws.Range("G10").Value = 60.43484133333334
$ws.Range("H10").Value = 181.304524
$ws.Range("I10").Value = 0.1450120099461104
$ws.Range("J10").Value = 0.1450120099461104
$ws.Range("M10").Value = 25.766648
$ws.Range("N10").Value = 77.299944
$ws.Range("O10").Value = 0.2921247832814228
$ws.Range("P10").Value = 0.2921247832814228
$ws.Range("Q10").Value = 1557.203283571851
$ws.Range("R10").Value = 14014.82955214666
$ws.Range("S10").Value = 0.04236160197871101
$ws.Range("T10").Value = 0.04236160197871102

# Row 11
$ws.Range("G11").Value = 89.46554166666668
$ws.Range("H11").Value = 268.396625
$ws.Range("I11").Value = 0.2146705068098712
$ws.Range("J11").Value = 0.2146705068098712
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 13.83830466666667
$ws.Range("N11").Value = 41.514914
$ws.Range("O11").Value = 0.1568893148900199
$ws.Range("P11").Value = 0.1568893148900199
$ws.Range("Q11").Value = 1238.051422751695
$ws.Range("R11").Value = 11142.46280476525
$ws.Range("S11").Value = 0.03367950874049404
$ws.Range("T11").Value = 0.03367950874049405

# Row 12
$ws.Range("G12").Value = 89.46554166666668
$ws.Range("H12").Value = 268.396625
$ws.Range("I12").Value = 0.2146705068098712
$ws.Range("J12").Value = 0.2146705068098712
$ws.Range("O12").Value = 0.5509859018285573
$ws.Range("P12").Value = 0.5509859018285573
$ws.Range("Q12").Value = 4347.962639477139
$ws.Range("R12").Value = 39131.66375529426
$ws.Range("S12").Value = 0.1182804227906303
$ws.Range("T12").Value = 0.1182804227906304

# Row 13
$ws.Range("G13").Value = 89.46554166666668
$ws.Range("H13").Value = 268.396625
$ws.Range("I13").Value = 0.2146705068098712
$ws.Range("J13").Value = 0.2146705068098712
$ws.Range("M13").Value = 25.766648
$ws.Range("N13").Value = 77.299944
$ws.Range("O13").Value = 0.2921247832814228
$ws.Range("P13").Value = 0.2921247832814228
$ws.Range("Q13").Value = 2305.227120254334
$ws.Range("R13").Value = 20747.044082289
$ws.Range("S13").Value = 0.06271057527874682
$ws.Range("T13").Value = 0.06271057527874682
